$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Mapping of (row, column) -> new value, in document order, matching the
# data rows of the 5-column practice table (rows 1, 5, 10, 15, 20).
$updates = @(
    @{ Row = 1;  Col = 1; Text = "140×6=840" },
    @{ Row = 1;  Col = 2; Text = "345×6=2070" },
    @{ Row = 1;  Col = 3; Text = "253×2=506" },
    @{ Row = 1;  Col = 4; Text = "394×2=788" },
    @{ Row = 1;  Col = 5; Text = "393×2=786" },

    @{ Row = 5;  Col = 1; Text = "895×7=6265" },
    @{ Row = 5;  Col = 2; Text = "522×4=2088" },
    @{ Row = 5;  Col = 3; Text = "690×3=2070" },
    @{ Row = 5;  Col = 4; Text = "353×2=706" },
    @{ Row = 5;  Col = 5; Text = "160×5=800" },

    @{ Row = 10; Col = 1; Text = "645×4=2580" },
    @{ Row = 10; Col = 2; Text = "906×7=6342" },
    @{ Row = 10; Col = 3; Text = "531×2=1062" },
    @{ Row = 10; Col = 4; Text = "487×4=1948" },
    @{ Row = 10; Col = 5; Text = "726×5=3630" },

    @{ Row = 15; Col = 1; Text = "318×5=1590" },
    @{ Row = 15; Col = 2; Text = "869×4=3476" },
    @{ Row = 15; Col = 3; Text = "216×9=1944" },
    @{ Row = 15; Col = 4; Text = "891×8=7128" },
    @{ Row = 15; Col = 5; Text = "376×3=1128" },

    @{ Row = 20; Col = 1; Text = "517×2=1034" },
    @{ Row = 20; Col = 2; Text = "358×2=716" },
    @{ Row = 20; Col = 3; Text = "537×3=1611" },
    @{ Row = 20; Col = 4; Text = "263×2=526" },
    @{ Row = 20; Col = 5; Text = "690×7=4830" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
